# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook

# --- Update F2:F9 timestamps on the "data" sheet ---
$data = $wb.Worksheets.Item("data")
$data.Range("F2").Value = "2021-10-05 14:22:41.485589"
$data.Range("F3").Value = "2021-10-05 14:22:41.485597"
$data.Range("F4").Value = "2021-10-05 14:22:41.485600"
$data.Range("F5").Value = "2021-10-05 14:22:41.485603"
$data.Range("F6").Value = "2021-10-05 14:22:41.485606"
$data.Range("F7").Value = "2021-10-05 14:22:41.485608"
$data.Range("F8").Value = "2021-10-05 14:22:41.485611"
$data.Range("F9").Value = "2021-10-05 14:22:41.485613"

# --- Add a new "metadata" worksheet positioned after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - reuse the exact header style ("data" sheet's B1 style)
# by copy/paste-special of formats, so the style index is shared rather than
# a brand-new style entry being minted.
$data.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# A2 reuses the "data" sheet's index-column style too.
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Severe multi-system atopic disease with high IgE"
$meta.Range("C2").Value = 62

# D2 must be the literal text "1.8" (not the number 1.8) with the sheet's
# default (no explicit) cell style. Building the string via TEXT() in a
# scratch cell and copying VALUES ONLY keeps it a genuine text cell without
# ever touching NumberFormat/quote-prefix on D2 (which would otherwise mint
# a brand-new, permanent style-table entry).
$meta.Range("Z1").Formula = '=TEXT(1.8,"0.0")'
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("E2").Value = "2021-07-13T14:30:23.111196Z"
$meta.Range("F2").Value = "2021-10-05 14:22:41.482294"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/62/?format=json"
